# Weekly data refresh: a new record (week of 2023-03-30) is inserted at the
# top of the "Feria Lagunitas de Puerto Montt - Apio" date-ordered block
# (row 360), pushing the existing rows 360:435 down to 361:436.
#
# Insert a new row at 360 - Excel shifts rows 360:435 -> 361:436 and
# extends the used range/dimension to A1:R436 automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(360).Insert()

# Populate the newly inserted row 360 with the new weekly record. Columns
# that are constant for this market/product/quality/unit carry over
# unchanged from the row that used to occupy 360 (now row 361).
$ws.Range("A360").Value = 4
$ws.Range("B360").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C360").Value = 'Los Lagos'
$ws.Range("D360").Value = 45015
$ws.Range("E360").Value = 10
$ws.Range("F360").Value = 100112017
$ws.Range("G360").Value = 'Apio'
$ws.Range("H360").Value = 'Americana (o)'
$ws.Range("I360").Value = 'Primera'
$ws.Range("J360").Value = 25
$ws.Range("K360").Value = 12000
$ws.Range("L360").Value = 12000
$ws.Range("M360").Value = 12000
$ws.Range("N360").Value = '$/docena de matas'
$ws.Range("O360").Value = 'Región de Coquimbo'
$ws.Range("P360").Value = 2000
$ws.Range("Q360").Value = 6
$ws.Range("R360").Value = 'Hortaliza'
